$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Safety Net" powerup row (row 24) to each of the six
# Powerup/Chance tables (columns A:B, E:F, I:J, M:N, Q:R, U:V), each with a
# 4% (0.04) chance.
$ws.Range("A24").Value = "Safety Net"
$ws.Range("B24").Value = 0.04

$ws.Range("E24").Value = "Safety Net"
$ws.Range("F24").Value = 0.04

$ws.Range("I24").Value = "Safety Net"
$ws.Range("J24").Value = 0.04

$ws.Range("M24").Value = "Safety Net"
$ws.Range("N24").Value = 0.04

$ws.Range("Q24").Value = "Safety Net"
$ws.Range("Q24").Style = "Normal"
$ws.Range("R24").Value = 0.04

$ws.Range("U24").Value = "Safety Net"
$ws.Range("U24").Style = "Normal"
$ws.Range("V24").Value = 0.04

# Extend the "Total %:" SUM formulas on row 25 to include the new row 24.
$ws.Range("B25").Formula = "=SUM(B5:B24)"
$ws.Range("F25").Formula = "=SUM(F5:F24)"
$ws.Range("J25").Formula = "=SUM(J5:J24)"
$ws.Range("N25").Formula = "=SUM(N5:N24)"
$ws.Range("R25").Formula = "=SUM(R5:R24)"
$ws.Range("V25").Formula = "=SUM(V5:V24)"

# Update the sheet view selection to the new last data cell.
$ws.Range("U24").Select()
